$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-14 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-15 Monday", 2) | Out-Null
$d.Content.Find.Execute("429÷2=214, 1", $true, $false, $false, $false, $false, $true, 1, $false, "779÷8=97, 3", 2) | Out-Null
$d.Content.Find.Execute("942÷2=471, 0", $true, $false, $false, $false, $false, $true, 1, $false, "520÷5=104, 0", 2) | Out-Null
$d.Content.Find.Execute("646÷8=80, 6", $true, $false, $false, $false, $false, $true, 1, $false, "573÷6=95, 3", 2) | Out-Null
$d.Content.Find.Execute("740÷5=148, 0", $true, $false, $false, $false, $false, $true, 1, $false, "402÷5=80, 2", 2) | Out-Null
$d.Content.Find.Execute("351÷5=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "400÷9=44, 4", 2) | Out-Null
$d.Content.Find.Execute("616÷3=205, 1", $true, $false, $false, $false, $false, $true, 1, $false, "762÷9=84, 6", 2) | Out-Null
$d.Content.Find.Execute("597÷3=199, 0", $true, $false, $false, $false, $false, $true, 1, $false, "843÷9=93, 6", 2) | Out-Null
$d.Content.Find.Execute("313÷2=156, 1", $true, $false, $false, $false, $false, $true, 1, $false, "479÷5=95, 4", 2) | Out-Null
$d.Content.Find.Execute("503÷3=167, 2", $true, $false, $false, $false, $false, $true, 1, $false, "304÷9=33, 7", 2) | Out-Null
$d.Content.Find.Execute("143÷2=71, 1", $true, $false, $false, $false, $false, $true, 1, $false, "699÷3=233, 0", 2) | Out-Null
$d.Content.Find.Execute("725÷2=362, 1", $true, $false, $false, $false, $false, $true, 1, $false, "871÷5=174, 1", 2) | Out-Null
$d.Content.Find.Execute("648÷5=129, 3", $true, $false, $false, $false, $false, $true, 1, $false, "673÷5=134, 3", 2) | Out-Null
$d.Content.Find.Execute("932÷4=233, 0", $true, $false, $false, $false, $false, $true, 1, $false, "847÷5=169, 2", 2) | Out-Null
$d.Content.Find.Execute("442÷9=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "374÷2=187, 0", 2) | Out-Null
$d.Content.Find.Execute("593÷8=74, 1", $true, $false, $false, $false, $false, $true, 1, $false, "605÷9=67, 2", 2) | Out-Null
$d.Content.Find.Execute("311÷5=62, 1", $true, $false, $false, $false, $false, $true, 1, $false, "756÷8=94, 4", 2) | Out-Null
$d.Content.Find.Execute("975÷3=325, 0", $true, $false, $false, $false, $false, $true, 1, $false, "260÷4=65, 0", 2) | Out-Null
$d.Content.Find.Execute("555÷2=277, 1", $true, $false, $false, $false, $false, $true, 1, $false, "588÷8=73, 4", 2) | Out-Null
$d.Content.Find.Execute("377÷6=62, 5", $true, $false, $false, $false, $false, $true, 1, $false, "179÷8=22, 3", 2) | Out-Null
$d.Content.Find.Execute("411÷5=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "576÷2=288, 0", 2) | Out-Null
$d.Content.Find.Execute("337÷9=37, 4", $true, $false, $false, $false, $false, $true, 1, $false, "749÷9=83, 2", 2) | Out-Null
$d.Content.Find.Execute("441÷3=147, 0", $true, $false, $false, $false, $false, $true, 1, $false, "586÷7=83, 5", 2) | Out-Null
$d.Content.Find.Execute("758÷5=151, 3", $true, $false, $false, $false, $false, $true, 1, $false, "564÷2=282, 0", 2) | Out-Null
$d.Content.Find.Execute("214÷9=23, 7", $true, $false, $false, $false, $false, $true, 1, $false, "132÷8=16, 4", 2) | Out-Null
$d.Content.Find.Execute("299÷6=49, 5", $true, $false, $false, $false, $false, $true, 1, $false, "453÷7=64, 5", 2) | Out-Null
